# Updated the meta data indicator
# Fills in column E ("meta_data") for rows 2-26 on Sheet1 with the
# corresponding SDG-indicator metadata string, matching the styling
# already used by the rest of each row, and restores the last active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row data: new "meta_data" values for E2:E26 -----------------------
$ws.Range("E2").Value  = "1.2.1, 10.1.1,10.2.1"
$ws.Range("E3").Value  = "3.8.1,3.8.2,1.a.2"
$ws.Range("E4").Value  = "2.1.2,2.2.1,2.2.2,2.2.3"
$ws.Range("E5").Value  = "3.1, 3.2, 3.8"
$ws.Range("E6").Value  = "3.5, 11.2"
$ws.Range("E7").Value  = "4.5, 5.b, 8.5, 8.6, 8.b, 9.2, 9.c"
$ws.Range("E8").Value  = "12.8,13.3"
$ws.Range("E9").Value  = "5.2.2, 11.7.2, 16.1.3, 16.2.3"
$ws.Range("E10").Value = "16.7.1"
$ws.Range("E11").Value = "6.3.2, 6.4.1, 6.4.2, 6.5.1, 6.5.2, 15.3.1"
$ws.Range("E12").Value = "3.9.1"
$ws.Range("E13").Value = "any economic statistics related SDG indicator"
$ws.Range("E14").Value = "12.2.1"
$ws.Range("E15").Value = "9.3.2"
$ws.Range("E16").Value = "1.1.1, 1.2.1, 10.2.1"
$ws.Range("E17").Value = "11.3.1, 11.7.1, 9.1.1"
$ws.Range("E18").Value = "11.2.1, 11.6.2, 11.7.1, 11.a.1, 15.1.2, 11.7.2, 11.b.1"
$ws.Range("E19").Value = "8.4.1"
$ws.Range("E20").Value = "4.7, 13.3"
$ws.Range("E21").Value = "4.7, 12.8"
$ws.Range("E22").Value = "14.7.1"
$ws.Range("E23").Value = "15.2.1"
$ws.Range("E24").Value = "14.5.1, 15.1.2"
$ws.Range("E25").Value = "16.1.1, 16.1.3, 16.1.4, 16.4.2, 16.a.1, 10.3.1, 16.b.1, 5.2.2"
$ws.Range("E26").Value = "17.1.2"

# --- Formatting --------------------------------------------------------
# Rows 16-26 already carry a distinct cell style (the "s=2" black-font
# style) on columns A-D; mirror that style onto the newly filled E cells
# by copying the formatting from the neighboring D column.
$ws.Range("D16:D26").Copy() | Out-Null
$ws.Range("E16:E26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Selection -----------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
